$wb = $excel.ActiveWorkbook
$wsSettings = $wb.Worksheets.Item("Settings")

# --- New "Transaction File (Input File)" section on the Settings sheet ---

# Section header row (row 8) with highlighted fill across A:C
$wsSettings.Range("A8").Value = "Transaction File (Input File)"
$wsSettings.Range("A8:C8").Interior.Color = 14282978

# Row 9 - Substitute Items sheet name setting
$wsSettings.Range("A9").Value = "Sheet_SubstituteItemsName"
$wsSettings.Range("B9").Value = "Substitute Items"
$wsSettings.Range("C9").Value = "Substitute Items Sheet name"

# Row 10 - Substitute Items file name setting
$wsSettings.Range("A10").Value = "File_SubstituteItemsName"
$wsSettings.Range("B10").Value = "Input - Substitute Items.xlsx"
$wsSettings.Range("C10").Value = "Substitute Items file name (Input file)"

# Row 11 - Buyers List sheet name setting
$wsSettings.Range("A11").Value = "Sheet_BuyersName"
$wsSettings.Range("B11").Value = "Buyers List"
$wsSettings.Range("C11").Value = "Buyers list sheet name"

# --- Active sheet / selection changes ---
# Settings becomes the active tab (previously it was Assets), with A9
# selected. Activating Settings also drops the stale "tabSelected" flag
# that was sitting on the Assets sheet's view.
[void]$wsSettings.Activate()
[void]$wsSettings.Range("A9").Select()
